$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2").Value = 10
$ws.Range("AH2").Value = 8
$ws.Range("BF2").Value = "'2013-03-21"
$ws.Range("AD3").Value = 21
$ws.Range("AF3").Value = 12
$ws.Range("AG3").Value = 14
$ws.Range("AI3").Value = 16
$ws.Range("BC3").Value = 14
$ws.Range("BF3").Value = "'2013-03-21"
$ws.Range("AD4").Value = 10
$ws.Range("AH4").Value = 8
$ws.Range("AI4").Value = 27
$ws.Range("AV4").Value = 14
$ws.Range("BF4").Value = "'2013-03-21"
$ws.Range("AD5").Value = 10
$ws.Range("AJ5").Value = 18
$ws.Range("BF5").Value = "'2013-03-21"
$ws.Range("D6").Value = 66
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = 0.545
$ws.Range("I6").Value = 35.5
$ws.Range("J6").Value = 81.59999999999999
$ws.Range("K6").Value = 0.435
$ws.Range("N6").Value = 0.34
$ws.Range("O6").Value = 16.8
$ws.Range("P6").Value = 21.3
$ws.Range("Q6").Value = 0.788
$ws.Range("U6").Value = 23
$ws.Range("V6").Value = 14.7
$ws.Range("Z6").Value = 19.6
$ws.Range("AA6").Value = 19.8
$ws.Range("AB6").Value = 92.8
$ws.Range("AC6").Value = 0.5
$ws.Range("AD6").Value = 28
$ws.Range("AF6").Value = 10
$ws.Range("AG6").Value = 12
$ws.Range("AJ6").Value = 15
$ws.Range("AU6").Value = 9
$ws.Range("AV6").Value = 15
$ws.Range("AW6").Value = 19
$ws.Range("BC6").Value = 13
$ws.Range("BF6").Value = "'2013-03-21"
$ws.Range("AD7").Value = 10
$ws.Range("AR7").Value = 7
$ws.Range("BF7").Value = "'2013-03-21"
$ws.Range("AD8").Value = 10
$ws.Range("AP8").Value = 23
$ws.Range("AU8").Value = 7
$ws.Range("BF8").Value = "'2013-03-21"
$ws.Range("D9").Value = 69
$ws.Range("E9").Value = 47
$ws.Range("G9").Value = 0.681
$ws.Range("J9").Value = 85.40000000000001
$ws.Range("L9").Value = 6.4
$ws.Range("M9").Value = 18.9
$ws.Range("N9").Value = 0.342
$ws.Range("Q9").Value = 0.695
$ws.Range("S9").Value = 31.8
$ws.Range("T9").Value = 45.2
$ws.Range("U9").Value = 24.4
$ws.Range("V9").Value = 15.2
$ws.Range("W9").Value = 9.199999999999999
$ws.Range("AB9").Value = 106.1
$ws.Range("AD9").Value = 2
$ws.Range("AN9").Value = 25
$ws.Range("AQ9").Value = 29
$ws.Range("BF9").Value = "'2013-03-21"
$ws.Range("AD10").Value = 2
$ws.Range("AQ10").Value = 28
$ws.Range("AW10").Value = 28
$ws.Range("BF10").Value = "'2013-03-21"
$ws.Range("AF11").Value = 12
$ws.Range("AW11").Value = 27
$ws.Range("AX11").Value = 25
$ws.Range("BF11").Value = "'2013-03-21"
$ws.Range("AD12").Value = 10
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 13
$ws.Range("AH12").Value = 24
$ws.Range("BF12").Value = "'2013-03-21"
$ws.Range("AD13").Value = 10
$ws.Range("BB13").Value = 20
$ws.Range("BF13").Value = "'2013-03-21"
$ws.Range("AD14").Value = 2
$ws.Range("AE14").Value = 4
$ws.Range("AG14").Value = 5
$ws.Range("BF14").Value = "'2013-03-21"
$ws.Range("AD15").Value = 2
$ws.Range("AN15").Value = 15
$ws.Range("AR15").Value = 14
$ws.Range("BF15").Value = "'2013-03-21"
$ws.Range("AD16").Value = 21
$ws.Range("AH16").Value = 16
$ws.Range("AK16").Value = 19
$ws.Range("AP16").Value = 24
$ws.Range("BF16").Value = "'2013-03-21"
$ws.Range("AD17").Value = 21
$ws.Range("BF17").Value = "'2013-03-21"
$ws.Range("AD18").Value = 21
$ws.Range("AH18").Value = 16
$ws.Range("BF18").Value = "'2013-03-21"
$ws.Range("D19").Value = 65
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 0.354
$ws.Range("I19").Value = 35.4
$ws.Range("J19").Value = 81.5
$ws.Range("L19").Value = 5.4
$ws.Range("M19").Value = 18
$ws.Range("N19").Value = 0.298
$ws.Range("O19").Value = 18.3
$ws.Range("P19").Value = 25.1
$ws.Range("Q19").Value = 0.73
$ws.Range("R19").Value = 12.4
$ws.Range("U19").Value = 21.8
$ws.Range("Z19").Value = 18.6
$ws.Range("AB19").Value = 94.40000000000001
$ws.Range("AD19").Value = 30
$ws.Range("AF19").Value = 21
$ws.Range("AG19").Value = 23
$ws.Range("AI19").Value = 28
$ws.Range("AJ19").Value = 19
$ws.Range("AR19").Value = 9
$ws.Range("BB19").Value = 22
$ws.Range("BF19").Value = "'2013-03-21"
$ws.Range("AD20").Value = 2
$ws.Range("AT20").Value = 21
$ws.Range("BF20").Value = "'2013-03-21"
$ws.Range("AD21").Value = 28
$ws.Range("AJ21").Value = 17
$ws.Range("BF21").Value = "'2013-03-21"
$ws.Range("AD22").Value = 2
$ws.Range("AH22").Value = 14
$ws.Range("BF22").Value = "'2013-03-21"
$ws.Range("AD23").Value = 2
$ws.Range("AJ23").Value = 8
$ws.Range("AU23").Value = 8
$ws.Range("BF23").Value = "'2013-03-21"
$ws.Range("D24").Value = 67
$ws.Range("F24").Value = 41
$ws.Range("G24").Value = 0.388
$ws.Range("I24").Value = 37.1
$ws.Range("J24").Value = 83.8
$ws.Range("K24").Value = 0.443
$ws.Range("N24").Value = 0.358
$ws.Range("P24").Value = 16.6
$ws.Range("Q24").Value = 0.719
$ws.Range("S24").Value = 30.6
$ws.Range("T24").Value = 41.3
$ws.Range("V24").Value = 13.1
$ws.Range("W24").Value = 7.3
$ws.Range("Y24").Value = 4.8
$ws.Range("Z24").Value = 18.7
$ws.Range("AA24").Value = 16.2
$ws.Range("AB24").Value = 92.3
$ws.Range("AC24").Value = -4
$ws.Range("AD24").Value = 21
$ws.Range("AK24").Value = 20
$ws.Range("AN24").Value = 16
$ws.Range("AS24").Value = 15
$ws.Range("AT24").Value = 20
$ws.Range("AW24").Value = 20
$ws.Range("AY24").Value = 12
$ws.Range("BF24").Value = "'2013-03-21"
$ws.Range("AD25").Value = 2
$ws.Range("BB25").Value = 21
$ws.Range("BF25").Value = "'2013-03-21"
$ws.Range("D26").Value = 67
$ws.Range("E26").Value = 31
$ws.Range("G26").Value = 0.463
$ws.Range("M26").Value = 23.6
$ws.Range("N26").Value = 0.351
$ws.Range("O26").Value = 16.3
$ws.Range("P26").Value = 20.9
$ws.Range("S26").Value = 30.4
$ws.Range("W26").Value = 6.9
$ws.Range("AA26").Value = 19.1
$ws.Range("AC26").Value = -1.4
$ws.Range("AD26").Value = 21
$ws.Range("AE26").Value = 19
$ws.Range("AG26").Value = 19
$ws.Range("AH26").Value = 6
$ws.Range("AO26").Value = 20
$ws.Range("AP26").Value = 22
$ws.Range("AS26").Value = 16
$ws.Range("AW26").Value = 26
$ws.Range("BF26").Value = "'2013-03-21"
$ws.Range("D27").Value = 68
$ws.Range("E27").Value = 24
$ws.Range("G27").Value = 0.353
$ws.Range("J27").Value = 83.7
$ws.Range("K27").Value = 0.445
$ws.Range("M27").Value = 19.9
$ws.Range("N27").Value = 0.369
$ws.Range("O27").Value = 17.8
$ws.Range("Q27").Value = 0.771
$ws.Range("R27").Value = 11.5
$ws.Range("S27").Value = 28.9
$ws.Range("Y27").Value = 6.3
$ws.Range("AC27").Value = -5.1
$ws.Range("AD27").Value = 10
$ws.Range("AG27").Value = 24
$ws.Range("AH27").Value = 12
$ws.Range("AJ27").Value = 7
$ws.Range("AR27").Value = 13
$ws.Range("AX27").Value = 26
$ws.Range("BF27").Value = "'2013-03-21"
$ws.Range("AD28").Value = 10
$ws.Range("AH28").Value = 8
$ws.Range("AV28").Value = 13
$ws.Range("BF28").Value = "'2013-03-21"
$ws.Range("AD29").Value = 10
$ws.Range("AF29").Value = 21
$ws.Range("AG29").Value = 21
$ws.Range("AN29").Value = 24
$ws.Range("BF29").Value = "'2013-03-21"
$ws.Range("AD30").Value = 10
$ws.Range("AH30").Value = 8
$ws.Range("AJ30").Value = 14
$ws.Range("AR30").Value = 7
$ws.Range("BF30").Value = "'2013-03-21"
$ws.Range("AD31").Value = 21
$ws.Range("AE31").Value = 22
$ws.Range("AF31").Value = 23
$ws.Range("AG31").Value = 22
$ws.Range("AJ31").Value = 16
$ws.Range("AT31").Value = 8
$ws.Range("BF31").Value = "'2013-03-21"
